$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Problem" text for rows 2-4
$ws.Range("F2").Value = "opa"
$ws.Range("F3").Value = "AKAKAKAKAKAKAKAKAKA"
$ws.Range("F4").Value = "OQQOQOAOAOA"

# Row 5's date moves from 23/12/2024 to 24/12/2024 (its other columns,
# including the Problem text, stay as they were)
$ws.Range("A5").Value = "24/12/2024"

# The old row 6 (a duplicate of the other "Não consigo acessar minha conta"
# rows) is removed entirely, shifting nothing else up since it was last
$ws.Rows.Item(6).EntireRow.Delete()
